$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = "MCT-3A-Microcontroladores"

$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "MCT-3A-Microcontroladores"

$ws.Range("E4").Value = "MCT-3A-Microcontroladores"

$ws.Range("C6").Value = "-"

$ws.Range("F7").Value = "MCT-3A-Microcontroladores"

$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = "MCT-3A-Microcontroladores"
